$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51
$ws.Range("F51").Value = 62
$ws.Range("G51").Value = 1017.42

# Row 57
$ws.Range("F57").Value = 1
$ws.Range("G57").Value = 59.49

# Row 61
$ws.Range("B61").Value = 25300.89

# Row 102
$ws.Range("F102").Value = 74
$ws.Range("G102").Value = 5248.08

# Row 107
$ws.Range("F107").Value = 44
$ws.Range("G107").Value = 3094.08

# Row 122
$ws.Range("F122").Value = 89
$ws.Range("G122").Value = 9118.940000000001

# Row 133
$ws.Range("B133").Value = 206303.75

# Row 181
$ws.Range("B181").Value = 57756

# Row 182
$ws.Range("B182").Value = 53925

# Row 192
$ws.Range("F192").Value = 4
$ws.Range("G192").Value = 10075

# Row 195
$ws.Range("B195").Value = 39219.11

# Row 221
$ws.Range("F221").Value = 52
$ws.Range("G221").Value = 3901.04

# Row 228
$ws.Range("B228").Value = 11584.57

# Row 314
$ws.Range("B314").Value = 61610
$ws.Range("D314").Value = 102.71
$ws.Range("E314").Value = 122.71
$ws.Range("F314").Value = 88
$ws.Range("G314").Value = 9038.48

# Row 315
$ws.Range("B315").Value = 57077
$ws.Range("D315").Value = 93.08
$ws.Range("E315").Value = 111.2
$ws.Range("F315").Value = 1
$ws.Range("G315").Value = 93.08

# Row 343
$ws.Range("F343").Value = 104
$ws.Range("G343").Value = 11643.84

# Row 380
$ws.Range("B380").Value = 253919.9

# Row 385
$ws.Range("F385").Value = 7
$ws.Range("G385").Value = 1423.17

# Row 389
$ws.Range("B389").Value = 22363.5

# Row 453
$ws.Range("F453").Value = 121
$ws.Range("G453").Value = 11888.25

# Row 466
$ws.Range("F466").Value = 45
$ws.Range("G466").Value = 2720.25

# Row 473
$ws.Range("B473").Value = 136136.34

# Row 491
$ws.Range("F491").Value = 468
$ws.Range("G491").Value = 6294.6

# Row 492
$ws.Range("F492").Value = 466
$ws.Range("G492").Value = 6127.9

# Row 493
$ws.Range("F493").Value = 544
$ws.Range("G493").Value = 6968.64

# Row 496
$ws.Range("F496").Value = 302
$ws.Range("G496").Value = 4961.86

# Row 499
$ws.Range("F499").Value = 294
$ws.Range("G499").Value = 1934.52

# Row 500
$ws.Range("F500").Value = 421
$ws.Range("G500").Value = 6828.62

# Row 501
$ws.Range("F501").Value = 85
$ws.Range("G501").Value = 1654.1

# Row 502
$ws.Range("F502").Value = 925
$ws.Range("G502").Value = 6086.5

# Row 505
$ws.Range("F505").Value = 407
$ws.Range("G505").Value = 5352.05

# Row 506
$ws.Range("F506").Value = 321
$ws.Range("G506").Value = 8442.299999999999

# Row 509
$ws.Range("B509").Value = 93826.17999999999

# Row 555
$ws.Range("F555").Value = 521
$ws.Range("G555").Value = 3542.8

# Row 556
$ws.Range("F556").Value = 350
$ws.Range("G556").Value = 2397.5

# Row 558
$ws.Range("F558").Value = 596
$ws.Range("G558").Value = 11830.6

# Row 559
$ws.Range("F559").Value = 318
$ws.Range("G559").Value = 2130.6

# Row 562
$ws.Range("F562").Value = 137
$ws.Range("G562").Value = 4471.68

# Row 563
$ws.Range("B563").Value = 36611.62

# Row 636
$ws.Range("F636").Value = 46
$ws.Range("G636").Value = 5650.18

# Row 640
$ws.Range("B640").Value = 208171.64

# Row 645
$ws.Range("F645").Value = 83
$ws.Range("G645").Value = 7880.85

# Row 646
$ws.Range("F646").Value = 7
$ws.Range("G646").Value = 190.4

# Row 647
$ws.Range("F647").Value = 5
$ws.Range("G647").Value = 136

# Row 649
$ws.Range("B649").Value = 53155.33

# Row 679
$ws.Range("F679").Value = 23
$ws.Range("G679").Value = 5164.19

# Row 682
$ws.Range("F682").Value = 11
$ws.Range("G682").Value = 902.4400000000001

# Row 684
$ws.Range("F684").Value = 22
$ws.Range("G684").Value = 6839.36

# Row 687
$ws.Range("F687").Value = 27
$ws.Range("G687").Value = 3147.66

# Row 689
$ws.Range("F689").Value = 25
$ws.Range("G689").Value = 2375

# Row 690
$ws.Range("F690").Value = 35
$ws.Range("G690").Value = 3702.65

# Row 693
$ws.Range("F693").Value = 1
$ws.Range("G693").Value = 99.31

# Row 695
$ws.Range("B695").Value = 37887.63

# Row 755
$ws.Range("F755").Value = 226
$ws.Range("G755").Value = 18432.56

# Row 758
$ws.Range("F758").Value = 252
$ws.Range("G758").Value = 32886

# Row 761
$ws.Range("F761").Value = 27
$ws.Range("G761").Value = 3011.58

# Row 763
$ws.Range("F763").Value = 99
$ws.Range("G763").Value = 2150.28

# Row 774
$ws.Range("F774").Value = 42
$ws.Range("G774").Value = 5069.82

# Row 775
$ws.Range("B775").Value = 245406.93

# Row 800
$ws.Range("F800").Value = 4
$ws.Range("G800").Value = 149.6

# Row 801
$ws.Range("B801").Value = 343.46

# Row 852
$ws.Range("F852").Value = 601
$ws.Range("G852").Value = 18168.23

# Row 853
$ws.Range("F853").Value = 3129
$ws.Range("G853").Value = 510371.19

# Row 856
$ws.Range("F856").Value = 117
$ws.Range("G856").Value = 4462.38

# Row 857
$ws.Range("F857").Value = 169
$ws.Range("G857").Value = 13036.66

# Row 858
$ws.Range("F858").Value = 113
$ws.Range("G858").Value = 16708.18

# Row 861
$ws.Range("B861").Value = 624482.9399999999

# Row 867
$ws.Range("B867").Value = 3478203.88

# Row 868
$ws.Range("B868").Value = 3478203.88
